# Update the three header-row labels (row 1) of the LandCover_Class table.
# "LCCC" -> "LAND COVER CODE", "DESCRIPTIO" -> "DESCRIPTION", "Walk_Impd" -> "IMPEDANCE"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "LAND COVER CODE"
$ws.Range("B1").Value = "DESCRIPTION"
$ws.Range("C1").Value = "IMPEDANCE"

# Give the new, longer A1 header its own look: centered, wrapped, same numeric
# display format as the rest of the header row, on a taller row.
$ws.Range("A1").NumberFormat = "0.000000"
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").WrapText = $true
$ws.Rows("1:1").RowHeight = 30

# Narrow column A now that the header wraps onto two lines.
$ws.Columns("A:A").ColumnWidth = 11.3

# Select the data table instead of the stray G13 cell left over from editing.
$ws.Range("A1:C30").Select() | Out-Null
